# Regenerate merged AHB files
#
# 1. Rename the header labels: the "_old" columns become the FV2404
#    (previous-format) columns and the "_new" columns become the FV2410
#    (current-format) columns.
# 2. Freeze the header row.
# 3. Turn the sheet's used range into a native Excel table ("Table1").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404",
    "diff",
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- Freeze the top (header) row -------------------------------------------
$ws.Activate()
[void]$ws.Range("A2").Select()
[void]($excel.ActiveWindow.FreezePanes = $true)

# --- Convert the used range into a table ------------------------------------
# Stash the header row's existing formatting (bold/fill/border/alignment) in
# an unused scratch row so it survives ListObjects.Add unaltered: adding a
# table otherwise snapshots any pre-existing header formatting that differs
# from the (new) table style into an explicit per-table dxf override, which
# the source workbook never had.
$headerRange = $ws.Range("A1:U1")
$scratch = $ws.Range("A100:U100")
$headerRange.Copy()
$scratch.PasteSpecial(-4122)

$headerRange.ClearFormats()

$dataRange = $ws.Range("A1:U57")
$table = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$table.Name = "Table1"
$table.TableStyle = $null

$scratch.Copy()
$headerRange.PasteSpecial(-4122)
$scratch.Clear()

[void]$ws.Range("A1").Select()

$wb.Save()
